$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the text of a *single* run range without ever going through
# Find/Replace (which silently runs AutoCorrect/AutoFormat -> smart quotes,
# smart apostrophes, etc.). Range.Text = "..." only ever rewrites the first
# run under the range and leaves the rest alone, so it is only safe to use
# when the range corresponds to exactly one run (e.g. located via Find, with
# MatchWholeWord / exact text so Start/End bracket only that run's text).
# ---------------------------------------------------------------------------

function Set-ExactText($rng, [string]$oldText, [string]$newText) {
    $f = $rng.Find
    $f.ClearFormatting()
    $ok = $f.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not locate text: $oldText"
    }
    $rng.Text = $newText
}

# 1) Title -------------------------------------------------------------
Set-ExactText $d.Content "Inviting Space: Challenges and Capabilities" "History's Echo: Lessons from the Annals of Time"

# 2) Author name ---------------------------------------------------------
Set-ExactText $d.Content "Elise Mayweather" "A.J. ""Sage"" Sinclair"

# 3) Email paragraph -------------------------------------------------------
Set-ExactText $d.Content "username@emaildomain" "historian"
Set-ExactText $d.Content "com" "ajsinclair@educonnect.org"

# ---------------------------------------------------------------------------
# 4) Main body paragraph (paragraph 5) - full rewrite.
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs(5)
$r = $p5.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Delete() | Out-Null

$ins = $d.Paragraphs(5).Range
$ins.Collapse(1) | Out-Null

$ins.InsertAfter("History, a sprawling tapestry of human experience, stands as an intricate, multi-dimensional story of civilizations, cultures, and countless lives interconnected across time's vast expanse.")
$ins.Collapse(0) | Out-Null
$ins.InsertAfter(" From the dawn of humanity, people have grappled with the mysteries of life, the challenges of existence, and the quest for meaning within the seemingly enigmatic tapestry of existence")
$ins.Collapse(0) | Out-Null
$ins.InsertAfter([char]11)
$ins.Collapse(0) | Out-Null
$ins.InsertAfter([char]11)
$ins.Collapse(0) | Out-Null
$ins.InsertAfter("Our ancestors, driven by insatiable curiosity, have left behind a treasure trove of knowledge and wisdom amassed over millennia.")
$ins.Collapse(0) | Out-Null
$ins.InsertAfter(" Etched in tomes, tablets, and inscribed artifacts, their narratives reveal civilizations birthed from humble origins, tracing their rise, zeniths, and inevitable ebbs.")
$ins.Collapse(0) | Out-Null
$ins.InsertAfter(" Delving into these chronicles transports us through time to learn from the trials and triumphs of civilizations past, gleaning insights into ourselves and our place in the grand narrative of humanity")
$ins.Collapse(0) | Out-Null
$ins.InsertAfter([char]11)
$ins.Collapse(0) | Out-Null
$ins.InsertAfter([char]11)
$ins.Collapse(0) | Out-Null
$ins.InsertAfter("History, woven with the aspirations, struggles, and resilience of generations, is not merely a collection of bygone events; it is a living testament to our capacity for both greatness and folly.")
$ins.Collapse(0) | Out-Null
$ins.InsertAfter(" Those who choose to delve into the fabric of time gain invaluable lessons, a deeper appreciation for the tapestry of life and an understanding of the interconnectedness of all humanity")
$ins.Collapse(0) | Out-Null
$ins.InsertAfter([char]11)
$ins.Collapse(0) | Out-Null
$ins.InsertAfter([char]11)
$ins.Collapse(0) | Out-Null
$ins.InsertAfter("Main Body:")
$ins.Collapse(0) | Out-Null
$ins.InsertAfter([char]11)
$ins.Collapse(0) | Out-Null
$ins.InsertAfter("1st Paragraph: Ancient wisdom, embedded in the annals of history, offers compelling lessons for contemporary societies.")
$ins.Collapse(0) | Out-Null
$ins.InsertAfter(" The rise and fall of empires, the echoes of diplomatic alliances or military conflicts, and the evolution of thought from antiquity to modernity provide fertile ground for learning and personal growth.")
$ins.Collapse(0) | Out-Null
$ins.InsertAfter(" By comprehending our origins and the triumphs and tragedies of civilizations past, we gain perspective on our current challenges and aspirations, leading to a more profound appreciation of the complexities of power, diplomacy, and governance.")
$ins.Collapse(0) | Out-Null
$ins.InsertAfter([char]11)
$ins.Collapse(0) | Out-Null
$ins.InsertAfter([char]11)
$ins.Collapse(0) | Out-Null
$ins.InsertAfter("2nd Paragraph: Historical research and analysis empower us to understand the ")
$ins.Collapse(0) | Out-Null
$ins.InsertAfter("roots of our present-day predicaments, be it deep-seated socio-political issues or complex cultural dynamics.")
$ins.Collapse(0) | Out-Null
$ins.InsertAfter(" With each page turned, we uncover patterns and insights that inform our choices, helping us navigate the intricacies of an interconnected world.")
$ins.Collapse(0) | Out-Null
$ins.InsertAfter(" We learn from past triumphs, grapple with the lessons of adversity, and grow in empathy, tolerance, and understanding of diverse perspectives.")
$ins.Collapse(0) | Out-Null
$ins.InsertAfter([char]11)
$ins.Collapse(0) | Out-Null
$ins.InsertAfter([char]11)
$ins.Collapse(0) | Out-Null
$ins.InsertAfter("3rd Paragraph: The study of history nurtures a sense of global citizenship, interconnectedness, and an appreciation for the contributions of diverse cultures throughout history.")
$ins.Collapse(0) | Out-Null
$ins.InsertAfter(" By recognizing the patterns woven through the ages, we cultivate an awareness of our shared humanity and learn to transcend boundaries and divisions.")
$ins.Collapse(0) | Out-Null
$ins.InsertAfter(" In this era of globalization, this is a priceless tool for building bridges across borders and cultivating global understanding and harmony")

$p5rng = $d.Paragraphs(5).Range
$p5rng.MoveEnd(1, -1) | Out-Null
$p5rng.Font.Name = "Aptos"
$p5rng.Font.Size = 12
$p5rng.Font.Color = 0

# ---------------------------------------------------------------------------
# 5) Summary body paragraph (paragraph 7) - full rewrite.
# ---------------------------------------------------------------------------
$p7 = $d.Paragraphs(7)
$r7 = $p7.Range
$r7.MoveEnd(1, -1) | Out-Null
$r7.Delete() | Out-Null

$ins7 = $d.Paragraphs(7).Range
$ins7.Collapse(1) | Out-Null
$ins7.InsertAfter("In exploring history, we navigate the annals of time, learning lessons from civilizations and cultures long gone.")
$ins7.Collapse(0) | Out-Null
$ins7.InsertAfter(" We gain insights into power, diplomacy, and governance, empathetically understanding diverse perspectives.")
$ins7.Collapse(0) | Out-Null
$ins7.InsertAfter(" Studying history allows us to grapple with current challenges, embrace global citizenship, and strive for a better future.")
$ins7.Collapse(0) | Out-Null
$ins7.InsertAfter(" Embracing the lessons it holds enables us to appreciate our place in the grand narrative of humanity and engage as responsible stewards of our shared legacy.")

$p7rng = $d.Paragraphs(7).Range
$p7rng.MoveEnd(1, -1) | Out-Null
$p7rng.Font.Name = "Aptos"
$p7rng.Font.Color = 0

# ---------------------------------------------------------------------------
# 6) Append a new trailing empty paragraph at the very end of the body.
# ---------------------------------------------------------------------------
$endRng = $d.Paragraphs(7).Range
$endRng.Collapse(0) | Out-Null
$endRng.InsertParagraphAfter()
